# V 2.0.2 se arreglo la fechar y hora de reimpresion
#
# Updates the patient-record form on the active sheet with the new
# patient's data (name, expediente, birth date, age, birth place,
# nationality, document number, emergency contact, and the date/time
# of the medical visit that was mis-printed on the previous reprint).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification -------------------------------------------------
$ws.Range("A6").Value = "RASCÒN  ALTALEF  JORGE  SAMUEL OSWALDO"
$ws.Range("G6").Value = "/201762650"

# --- Birth data --------------------------------------------------------------
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1973-06-08"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44"

$ws.Range("E9").Value = "TIQUISATE, ESCUINTLA"

# --- Ocupacion / Nacionalidad / Documento de identificacion ------------------
$ws.Range("C11").Value = ""
$ws.Range("E11").Value = "GUATEMALTECO"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "1878491610506"

# --- Emergency contact ---------------------------------------------------
$ws.Range("A13").Value = "SINTI RAXTUN"
$ws.Range("D13").Value = "HERMANA"
$ws.Range("E13").Value = "13 C. 3-08 Z. 2 EL ZAPOTE"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "52063422"

# --- Fecha / hora de la asistencia medica ------------------------------------
$ws.Range("D14").Value = "Hora: 15:46:8"
$ws.Range("E14").Value = "Area de urgencia: CIRUGIA"
$ws.Range("A15").Value = "24/10/2017"
